# Update the dSF (column F) values per the repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = -10
    6  = 8
    7  = 1
    8  = -3
    9  = 2
    10 = -5
    11 = -5
    13 = -3
    14 = 2
    15 = -1
    16 = -7
    17 = 2
    20 = -1
    21 = 5
    22 = -2
    23 = 2
    24 = 2
    25 = 3
    26 = 1
    28 = -5
    29 = -3
    31 = 0
    32 = 8
    33 = -1
    35 = -2
    36 = 3
    37 = 3
    38 = -3
    39 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
